$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3, shifting existing rows 3-15 down to 4-16
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new episode data
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Mon, 6 Dec 2021 10:55:00 +0000"
$ws.Cells.Item(3, 3).Value = "The Trial of Ghislaine Maxwell"
$ws.Cells.Item(3, 4).Value = "00:32:59"
$ws.Cells.Item(3, 5).Value = "https://dts.podtrac.com/redirect.mp3/chrt.fm/track/8DB4DB/pdst.fm/e/nyt.simplecastaudio.com/03d8b493-87fc-4bd1-931f-8a8e9b945d8a/episodes/4f9a4390-d07f-44dc-aabf-f87dece9fbd1/audio/128/default.mp3?aid=rss_feed&awCollectionId=03d8b493-87fc-4bd1-931f-8a8e9b945d8a&awEpisodeId=4f9a4390-d07f-44dc-aabf-f87dece9fbd1&feed=54nAGcIl"

# Match the bordered/centered style used by column A in the other data rows
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)
